$d = $word.ActiveDocument

# Locate the paragraph that contains the known sentence and append a new
# sentence to it as a new run (matching the existing eastAsia-hinted font).
$target = $d.Content.Find.Execute("晴，今天是高考第一天，上午考语文，下午考数学。")

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*晴，今天是高考第一天，上午考语文，下午考数学。*") {
        $r = $p.Range
        $r.Collapse(0)  # wdCollapseEnd
        $r.InsertAfter("今天天气不错")
        break
    }
}
